$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "DEPARTAMENTO"
$ws.Range("B1").Value = "LATITUD"
$ws.Range("C1").Value = "LONGITUD"
$ws.Range("D1").Value = "COLOR"

# Row 2 - AYACUCHO (unchanged values, kept for completeness)
$ws.Range("A2").Value = "AYACUCHO"
$ws.Range("B2").Value = -13.1581
$ws.Range("C2").Value = -74.2239
$ws.Range("D2").Value = "green"

# Row 3 - HUANCAVELICA (unchanged values, kept for completeness)
$ws.Range("A3").Value = "HUANCAVELICA"
$ws.Range("B3").Value = -12.7875
$ws.Range("C3").Value = -74.975
$ws.Range("D3").Value = "red"

# Row 4 - HUANUCO (name without accent, color changed to pink)
$ws.Range("A4").Value = "HUANUCO"
$ws.Range("B4").Value = -9.9306
$ws.Range("C4").Value = -76.2422

# Row 5 - JUNIN (unchanged values, kept for completeness)
$ws.Range("A5").Value = "JUNIN"
$ws.Range("B5").Value = -11.1574
$ws.Range("C5").Value = -75.9941
$ws.Range("D5").Value = "blue"

# Row 6 - PASCO (color changed to black)
$ws.Range("A6").Value = "PASCO"
$ws.Range("B6").Value = -10.6837
$ws.Range("C6").Value = -76.2567
$ws.Range("D6").Value = "black"

# Row 7 - SAN MARTIN (replaces UCAYALI, new coordinates, color purple)
$ws.Range("A7").Value = "SAN MARTIN"
$ws.Range("B7").Value = -6.51389
$ws.Range("C7").Value = -76.7408
$ws.Range("D7").Value = "purple"

# Set D4 after A7 so shared-string insertion order matches the target file
$ws.Range("D4").Value = "pink"

# Update the selection to match the target state
$ws.Range("E9").Select() | Out-Null
